# The underlying commit for this workbook ("Added a new file that I'm
# messing with for UI") does not correspond to any real data, formula, or
# formatting change inside this particular spreadsheet: every hunk in the
# recorded OOXML diff is a byte-level / serializer artifact (redundant
# xmlns declarations collapsing, an explicit-but-empty <numFmts count="0"/>,
# col width values being printed as "10.0" instead of "10", customWidth
# flags printing as "true" instead of "1", and every cell gaining an
# explicit s="0" default-style index) produced by the tool that re-saved
# the workbook outside of Excel. None of it changes any cell value,
# formula, shared string, row/column layout, or cell formatting that
# Excel's object model exposes.
#
# So the faithful COM replay of this commit is simply to open the
# workbook, touch it, and let Excel recalculate/resave it as-is -
# mirroring a user who only clicked around ("messing with UI") without
# actually changing any sheet content.

$wb = $excel.ActiveWorkbook

# Touch each existing worksheet (select it, force a recalculation) without
# writing any values/formulas/styles - this matches the no-op nature of
# the recorded diff while still exercising the COM surface.
foreach ($ws in $wb.Worksheets) {
    $ws.Calculate()
}

$wb.Worksheets.Item("Match Data").Activate()
$wb.Save()
